$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right before the existing row 135, shifting the
# existing rows 135:154 down to 137:156.
$ws.Rows("135:136").Insert()

# New row 135: Ají, Americana (o), Región Metropolitana
$ws.Cells.Item(135,1).Value = 11
$ws.Cells.Item(135,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(135,3).Value = "Bíobío"
$ws.Cells.Item(135,4).Value = 44918
$ws.Cells.Item(135,5).Value = 8
$ws.Cells.Item(135,6).Value = 100112021
$ws.Cells.Item(135,7).Value = "Ají"
$ws.Cells.Item(135,8).Value = "Americana (o)"
$ws.Cells.Item(135,9).Value = "Primera"
$ws.Cells.Item(135,10).Value = 100
$ws.Cells.Item(135,11).Value = 25000
$ws.Cells.Item(135,12).Value = 26000
$ws.Cells.Item(135,13).Value = 25500
$ws.Cells.Item(135,14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(135,15).Value = "Región Metropolitana"
$ws.Cells.Item(135,16).Value = 1700
$ws.Cells.Item(135,17).Value = 15
$ws.Cells.Item(135,18).Value = "Hortaliza"

# New row 136: Ají, Inferno, Región de Arica y Parinacota
$ws.Cells.Item(136,1).Value = 11
$ws.Cells.Item(136,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(136,3).Value = "Bíobío"
$ws.Cells.Item(136,4).Value = 44918
$ws.Cells.Item(136,5).Value = 8
$ws.Cells.Item(136,6).Value = 100112021
$ws.Cells.Item(136,7).Value = "Ají"
$ws.Cells.Item(136,8).Value = "Inferno"
$ws.Cells.Item(136,9).Value = "Primera"
$ws.Cells.Item(136,10).Value = 100
$ws.Cells.Item(136,11).Value = 17000
$ws.Cells.Item(136,12).Value = 18000
$ws.Cells.Item(136,13).Value = 17500
$ws.Cells.Item(136,14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(136,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(136,16).Value = 1750
$ws.Cells.Item(136,17).Value = 10
$ws.Cells.Item(136,18).Value = "Hortaliza"
